# Batch runner is added
#
# On the "SuiteDetails" sheet, re-tag the RevokeContract-related test rows
# (rows 52-143) in the "SuiteType" column (I) from "SmokeSuite" to the new
# "Regression" suite so they get picked up by the new batch/regression
# runner. A handful of rows in that span are section-header rows and have
# no SuiteType value, so they are skipped (90, 109, 120, 135, 137) -
# matching the existing conditional-formatting / data-validation block
# boundaries already on the sheet (J52:J89, J91:J108, J110:J119,
# J121:J134, J136, J138:J143).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SuiteDetails")

$ws.Range("I52:I89").Value = "Regression"
$ws.Range("I91:I108").Value = "Regression"
$ws.Range("I110:I119").Value = "Regression"
$ws.Range("I121:I134").Value = "Regression"
$ws.Range("I136:I136").Value = "Regression"
$ws.Range("I138:I143").Value = "Regression"

# Leave the sheet's cursor where the author left it while adding the rows.
[void]$ws.Activate()
$ws.Range("J146").Select() | Out-Null
